# Modification du document Réponses
# Adds the "Alors, on veut Ud à 0..." sentence (4 runs with distinct
# formatting) at the end of the 3rd bullet paragraph of the single
# textbox on slide 1 ("Que représentent les commandes Ud et Uq ...").

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Locate the insertion point: right after the trailing space that
# closes the "... control the torque through i_q.) " run, i.e. the
# very end of that bullet's paragraph (before the next bullet starts).
$full   = $tr.Text
$marker = "around 0 and control the torque through i_q.) "
$mIdx   = $full.IndexOf($marker)
if ($mIdx -lt 0) {
    throw "Could not locate insertion marker in shape text"
}
$anchorPos1 = $mIdx + $marker.Length   # 1-based Characters() position just after the space
$anchor = $tr.Characters($anchorPos1, 1)

# The 4 new runs to append (with their own run-level formatting).
# Inserted in reverse order against the same fixed anchor, since
# TextRange.InsertAfter always inserts immediately after the anchor
# (pushing previously-inserted text further along).
$runTexts = @(
    " Alors, on veut Ud à 0 car c’est une composante inutile pour la mise en rotation du moteur ; ainsi, si on transfert l’E utilisée sur Ud pour ",
    "Uq",
    ", on a un ",
    "fonctionnement optimal/maximal."
)

for ($i = $runTexts.Length - 1; $i -ge 0; $i--) {
    $anchor.InsertAfter($runTexts[$i]) | Out-Null
}

# Re-select each freshly inserted run (by character offset) and apply
# its italic formatting to force it onto its own <a:r> run.
$pos = $anchorPos1 + 1
foreach ($t in $runTexts) {
    $len = $t.Length
    $rng = $tr.Characters($pos, $len)
    $rng.Font.Italic = $true
    $pos += $len
}
